$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("scenarios")

# Insert a new column before column C (shifts periods..rich_win_prob right by one)
$ws.Columns("C:C").Insert()

# Populate the newly inserted column C with header and data
$ws.Range("C1").Value = "number_of_run"
$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 2
$ws.Range("C4").Value = 1

# Reflect the new selection state (entire column C selected, as after an insert-column action)
$ws.Range("C1:C1048576").Select()
